$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: voltage header for the new (10V, both pumps) test block
$ws.Range("A12").Value = "Tensão nas bombas"
$ws.Range("A12").HorizontalAlignment = -4108
$ws.Range("B12").Value = 5.5
$ws.Range("B12").HorizontalAlignment = -4108

# Row 13: column headers (bold + centered), reusing existing shared strings
$ws.Range("A13").Value = "Tanques"
$ws.Range("B13").Value = "Nível regua"
$ws.Range("C13").Value = "Nível SW"
$ws.Range("E13").Value = "Gamma 1"
$ws.Range("F13").Value = "Gamma 2"
$ws.Range("G13").Value = "Gamma1+Gamma2"
foreach ($col13 in @("A13","B13","C13","E13","F13","G13")) {
    $ws.Range($col13).HorizontalAlignment = -4108
    $ws.Range($col13).Font.Bold = $true
}

# Row 14: tank 1 data + gamma formulas
$ws.Range("A14").Value = 1
$ws.Range("A14").HorizontalAlignment = -4108
$ws.Range("B14").Value = 16.2
$ws.Range("B14").HorizontalAlignment = -4108
$ws.Range("C14").Value = 16
$ws.Range("C14").HorizontalAlignment = -4108
$ws.Range("E14").Formula = "=B14/(B14+B17)"
$ws.Range("F14").Formula = "=B15/(B15+B16)"
$ws.Range("G14").Formula = "=E14+F14"

# Row 15: tank 2 data
$ws.Range("A15").Value = 2
$ws.Range("A15").HorizontalAlignment = -4108
$ws.Range("B15").Value = 15.7
$ws.Range("B15").HorizontalAlignment = -4108
$ws.Range("C15").Value = 15.5
$ws.Range("C15").HorizontalAlignment = -4108

# Row 16: tank 3 data
$ws.Range("A16").Value = 3
$ws.Range("A16").HorizontalAlignment = -4108
$ws.Range("B16").Value = 10.7
$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("C16").Value = 6.6
$ws.Range("C16").HorizontalAlignment = -4108

# Row 17: tank 4 data
$ws.Range("A17").Value = 4
$ws.Range("A17").HorizontalAlignment = -4108
$ws.Range("B17").Value = 9.9
$ws.Range("B17").HorizontalAlignment = -4108
$ws.Range("C17").Value = 6.8
$ws.Range("C17").HorizontalAlignment = -4108

$ws.Range("C20").Select()
